$d = $word.ActiveDocument

# Italicize "St. Bernard"
$r1 = $d.Content
$r1.Find.Execute("St. Bernard", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r1.Italic = 1

# Italicize "collie"
$r2 = $d.Content
$r2.Find.Execute("collie", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r2.Italic = 1

# Bold "Presbyterian."
$r3 = $d.Content
$r3.Find.Execute("Presbyterian.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r3.Bold = 1
